# Add example values to the "New Table" worksheet (row 2 of the annotation
# table). Blank/unset cells keep their empty string; only cells that gain a
# concrete example value in the target revision are written here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Table")

$ws.Range("C2").Value  = "GO"
$ws.Range("D2").Value  = "https://bioregistry.io/GO:0019953"
$ws.Range("E2").Value  = "Fluvisol"
$ws.Range("H2").Value  = "'5.2"
$ws.Range("K2").Value  = "field"
$ws.Range("N2").Value  = "Diseased"
$ws.Range("O2").Value  = "NCIT"
$ws.Range("P2").Value  = "https://bioregistry.io/NCIT:C175249"
$ws.Range("Q2").Value  = "Cercospora zeae-maydis"
$ws.Range("T2").Value  = "diploid"
$ws.Range("U2").Value  = "PATO"
$ws.Range("V2").Value  = "https://bioregistry.io/PATO:0001394"
$ws.Range("W2").Value  = "'10"
$ws.Range("Z2").Value  = "2.4 Gb"
$ws.Range("AC2").Value = "maize"
$ws.Range("AF2").Value = "cultivar"
$ws.Range("AG2").Value = "EFO"
$ws.Range("AH2").Value = "https://bioregistry.io/EFO:0005136"
$ws.Range("AI2").Value = "B73"
$ws.Range("AL2").Value = "10 plants per 2 m row"
$ws.Range("AO2").Value = "soil"
$ws.Range("AP2").Value = "ENVO"
$ws.Range("AQ2").Value = "https://bioregistry.io/ENVO:00001998"
$ws.Range("AR2").Value = "KH2PO4 (170mg/L)"
$ws.Range("AU2").Value = "H3BO3 (6.2mg/L)"
$ws.Range("AX2").Value = "Nicotinic acid (0.5mg/L)"
$ws.Range("BA2").Value = "sucrose"
$ws.Range("BB2").Value = "CHEBI"
$ws.Range("BC2").Value = "https://bioregistry.io/CHEBI:17992"
$ws.Range("BD2").Value = "0.5mg/L NAA"
$ws.Range("BG2").Value = "agar"
$ws.Range("BH2").Value = "CHEBI"
$ws.Range("BI2").Value = "https://bioregistry.io/CHEBI:2509"
$ws.Range("BJ2").Value = "'5.2"
$ws.Range("BM2").Value = "25 °C day, 18 °C night"
$ws.Range("BP2").Value = "40 kg P2O5 per hectar"
$ws.Range("BS2").Value = "rain-fed"
$ws.Range("BV2").Value = "drought environment"
$ws.Range("BW2").Value = "EO"
$ws.Range("BX2").Value = "http://purl.obolibrary.org/obo/EO_0007404"
$ws.Range("BY2").Value = "1000 µmol m-2 s-1, 16 h"
$ws.Range("CB2").Value = "sample1"
